$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -------------------------------------------------------
# Row 2: the in-flight "Occupied Command" now targets a different room
# (RM-2109 instead of RM-4132), was issued at a later timestamp, and its
# set_value flips from -1 to 1.
$ws.Range("C2").Value = "RM-2109"
$ws.Range("E2").Value = "11/10/2015 14:35:15"
$ws.Range("F2").Value = 1

# Row 3 (new): a follow-up command for the same zone/actuator, issued a
# little later, restoring the original (-1) set_value.
$ws.Range("C3").Value = "RM-2109"
$ws.Range("D3").Value = "Occupied Command"
$ws.Range("E3").Value = "11/10/2015 14:50:15"
$ws.Range("F3").Value = -1

# --- Column widths ------------------------------------------------------
# Column C widened, column D given an explicit width for the first time.
$ws.Columns.Item(3).ColumnWidth = 14.666666666666666
$ws.Columns.Item(4).ColumnWidth = 15.166666666666666

# --- Selection ------------------------------------------------------
$ws.Range("D3").Select() | Out-Null
